$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N3").Value = 0.1
$ws.Range("O3").Value = 0.1

$ws.Range("O3").Select()
